# Auto update Excel log
# Appends new sensor-log rows to the PIR, Proximity and Camera sheets,
# mirroring the data the logging service captured for 2026-02-01.

$wb = $excel.ActiveWorkbook

# NOTE: this PowerShell-subset interpreter does not reliably bind
# PowerShell-style named arguments (`-paramName value`) on custom function
# calls, so the helper below is invoked with plain positional arguments.
function Append-LogRows($ws, $startRow, $rows) {
    $r = $startRow
    foreach ($row in $rows) {
        # Column A holds a YYYY-MM-DD date-looking string that must stay a
        # literal text value (not get auto-converted to a date serial), so
        # force the cell to Text format before assigning it.
        $cellA = $ws.Cells.Item($r, 1)
        $cellA.NumberFormat = "@"
        $cellA.Value = $row[0]

        $ws.Cells.Item($r, 2).Value = $row[1]
        $ws.Cells.Item($r, 3).Value = $row[2]
        $ws.Cells.Item($r, 4).Value = $row[3]
        $ws.Cells.Item($r, 5).Value = $row[4]
        $ws.Cells.Item($r, 6).Value = $row[5]

        $r = $r + 1
    }
}

# --- PIR sheet: 15 new bathroom motion readings (rows 2-16) ---
$pirSheet = $wb.Worksheets.Item("PIR")
$pirData = @(
    @("2026-02-01","13:42:56","13:00","Bathroom","No Motion","Inactive"),
    @("2026-02-01","13:43:01","13:00","Bathroom","No Motion","Inactive"),
    @("2026-02-01","13:43:06","13:00","Bathroom","No Motion","Inactive"),
    @("2026-02-01","13:43:11","13:00","Bathroom","No Motion","Inactive"),
    @("2026-02-01","13:43:16","13:00","Bathroom","No Motion","Inactive"),
    @("2026-02-01","13:43:20","13:00","Bathroom","No Motion","Inactive"),
    @("2026-02-01","13:43:21","13:00","Bathroom","No Motion","Inactive"),
    @("2026-02-01","13:43:25","13:00","Bathroom","No Motion","Inactive"),
    @("2026-02-01","13:43:26","13:00","Bathroom","No Motion","Inactive"),
    @("2026-02-01","13:43:30","13:00","Bathroom","No Motion","Inactive"),
    @("2026-02-01","13:43:31","13:00","Bathroom","No Motion","Inactive"),
    @("2026-02-01","13:43:36","13:00","Bathroom","No Motion","Inactive"),
    @("2026-02-01","13:43:37","13:00","Bathroom","Motion Detected","Active"),
    @("2026-02-01","13:43:45","13:00","Bathroom","No Motion","Inactive"),
    @("2026-02-01","13:43:50","13:00","Bathroom","No Motion","Inactive")
)
Append-LogRows $pirSheet 2 $pirData

# --- Proximity sheet: 2 new Living Room Main Door events (rows 18-19) ---
$proximitySheet = $wb.Worksheets.Item("Proximity")
$proximityData = @(
    @("2026-02-01","13:40:40","13:00","Living Room Main Door","ENTER","User ENTERED Living Room Main Door"),
    @("2026-02-01","13:40:40","13:00","Living Room Main Door","EXIT","User EXITED Living Room Main Door")
)
Append-LogRows $proximitySheet 18 $proximityData

# --- Camera sheet: 4 new captures for Living Room Main Door (rows 9-12) ---
$cameraSheet = $wb.Worksheets.Item("Camera")
$cameraData = @(
    @("2026-02-01","13:40:40","13:00","Living Room Main Door","Image Captured","Active"),
    @("2026-02-01","13:40:40","13:00","Living Room Main Door","Image Captured","Active"),
    @("2026-02-01","13:40:40","13:00","Living Room Main Door","Image Captured","Active"),
    @("2026-02-01","13:40:41","13:00","Living Room Main Door","Image Captured","Active")
)
Append-LogRows $cameraSheet 9 $cameraData
